$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 485.56
$ws.Range("J17").Value = 292.33334
$ws.Range("L17").Value = 877.0000200000001
$ws.Range("N17").Value = -1213.00002
$ws.Range("H106").Value = 2899
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 2899
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 2899
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -4161
$ws.Range("H116").Value = 362743.25
$ws.Range("I116").Value = 771461.5600000001
$ws.Range("J116").Value = 8520.733
$ws.Range("K116").Value = 771461.5600000001
$ws.Range("L116").Value = 8520.733
$ws.Range("M116").Value = -768019.5600000001
$ws.Range("N116").Value = -15404.733
$ws.Range("H132").Value = 253920.45
$ws.Range("I132").Value = 3731.8948
$ws.Range("K132").Value = 11195.6844
$ws.Range("M132").Value = -8665.6844
$ws.Range("H138").Value = 2782.9124
$ws.Range("I138").Value = 1962.5834
$ws.Range("J138").Value = 3001.6667
$ws.Range("K138").Value = 5887.7502
$ws.Range("L138").Value = 9005.000100000001
$ws.Range("M138").Value = -747.7502000000004
$ws.Range("N138").Value = -19285.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 8666.333000000001
$ws.Range("J19").Value = 8666.333000000001
$ws.Range("L19").Value = 8666.333000000001
$ws.Range("N19").Value = -9124.333000000001
$ws.Range("H61").Value = 1736.4445
$ws.Range("I61").Value = 1736.4445
$ws.Range("K61").Value = 1736.4445
$ws.Range("M61").Value = -1524.4445
$ws.Range("H74").Value = 2358.7778
$ws.Range("I74").Value = 1431.8667
$ws.Range("J74").Value = 6993.3335
$ws.Range("K74").Value = 1431.8667
$ws.Range("L74").Value = 6993.3335
$ws.Range("M74").Value = -557.8667
$ws.Range("N74").Value = -8741.333500000001
$ws.Range("H77").Value = 2358.7778
$ws.Range("I77").Value = 1431.8667
$ws.Range("J77").Value = 6993.3335
$ws.Range("K77").Value = 7159.333500000001
$ws.Range("L77").Value = 34966.6675
$ws.Range("M77").Value = -2791.333500000001
$ws.Range("N77").Value = -43702.6675
$ws.Range("H97").Value = 1001.6667
$ws.Range("I97").Value = 815.36365
$ws.Range("J97").Value = 1294.4286
$ws.Range("K97").Value = 815.36365
$ws.Range("L97").Value = 1294.4286
$ws.Range("M97").Value = -319.36365
$ws.Range("N97").Value = -2286.4286
$ws.Range("H136").Value = 1736.4445
$ws.Range("I136").Value = 1736.4445
$ws.Range("K136").Value = 5209.333500000001
$ws.Range("M136").Value = -2659.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7054.357
$ws.Range("I20").Value = 1576.1428
$ws.Range("J20").Value = 12532.571
$ws.Range("K20").Value = 1576.1428
$ws.Range("L20").Value = 12532.571
$ws.Range("M20").Value = -1329.1428
$ws.Range("N20").Value = -13026.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H15").Value = 14049.3
$ws.Range("I15").Value = 10008
$ws.Range("J15").Value = 14498.333
$ws.Range("K15").Value = 10008
$ws.Range("L15").Value = 14498.333
$ws.Range("M15").Value = -9838
$ws.Range("N15").Value = -14838.333
$ws.Range("H31").Value = 5572.649
$ws.Range("I31").Value = 2143.0605
$ws.Range("J31").Value = 10288.333
$ws.Range("K31").Value = 2143.0605
$ws.Range("L31").Value = 10288.333
$ws.Range("M31").Value = -1848.0605
$ws.Range("N31").Value = -10878.333
$ws.Range("H34").Value = 5572.649
$ws.Range("I34").Value = 2143.0605
$ws.Range("J34").Value = 10288.333
$ws.Range("K34").Value = 2143.0605
$ws.Range("L34").Value = 10288.333
$ws.Range("M34").Value = -1941.0605
$ws.Range("N34").Value = -10692.333
$ws.Range("H38").Value = 19993.5
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 19993.5
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 19993.5
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -20747.5
$ws.Range("H46").Value = 19993.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 19993.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 19993.5
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -20415.5
$ws.Range("H58").Value = 2215.1633
$ws.Range("I58").Value = 1404.1428
$ws.Range("J58").Value = 4242.7144
$ws.Range("K58").Value = 1404.1428
$ws.Range("L58").Value = 4242.7144
$ws.Range("M58").Value = -1201.1428
$ws.Range("N58").Value = -4648.7144
$ws.Range("H124").Value = 47998
$ws.Range("J124").Value = 47998
$ws.Range("L124").Value = 47998
$ws.Range("N124").Value = -52908
$ws.Range("H132").Value = 2603.04
$ws.Range("I132").Value = 1948.9546
$ws.Range("J132").Value = 7399.6665
$ws.Range("K132").Value = 5846.8638
$ws.Range("L132").Value = 22198.9995
$ws.Range("M132").Value = -3316.8638
$ws.Range("N132").Value = -27258.9995
$ws.Range("H136").Value = 2215.1633
$ws.Range("I136").Value = 1404.1428
$ws.Range("J136").Value = 4242.7144
$ws.Range("K136").Value = 4212.428400000001
$ws.Range("L136").Value = 12728.1432
$ws.Range("M136").Value = -1662.428400000001
$ws.Range("N136").Value = -17828.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 743402.1
$ws.Range("I5").Value = 660
$ws.Range("J5").Value = 1029072.1
$ws.Range("K5").Value = 1980
$ws.Range("L5").Value = 3087216.3
$ws.Range("M5").Value = -1868
$ws.Range("N5").Value = -3087440.3
$ws.Range("H103").Value = 816.6667
$ws.Range("I103").Value = 380
$ws.Range("J103").Value = 3000
$ws.Range("K103").Value = 1140
$ws.Range("L103").Value = 9000
$ws.Range("M103").Value = -261
$ws.Range("N103").Value = -10758
$ws.Range("H106").Value = 3995
$ws.Range("J106").Value = 3995
$ws.Range("L106").Value = 11985
$ws.Range("N106").Value = -13877
$ws.Range("H109").Value = 1407.8667
$ws.Range("I109").Value = 651.2857
$ws.Range("J109").Value = 12000
$ws.Range("K109").Value = 1953.8571
$ws.Range("L109").Value = 36000
$ws.Range("M109").Value = -913.8571000000002
$ws.Range("N109").Value = -38080
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H113").Value = 5682396
$ws.Range("I113").Value = 648.7778
$ws.Range("K113").Value = 1946.3334
$ws.Range("M113").Value = 223.6666
$ws.Range("H117").Value = 1555
$ws.Range("I117").Value = 400
$ws.Range("J117").Value = 1843.75
$ws.Range("K117").Value = 1200
$ws.Range("L117").Value = 5531.25
$ws.Range("M117").Value = 2242
$ws.Range("N117").Value = -12415.25
$ws.Range("H118").Value = 810.44446
$ws.Range("I118").Value = 549.25
$ws.Range("J118").Value = 2900
$ws.Range("K118").Value = 1647.75
$ws.Range("L118").Value = 8700
$ws.Range("M118").Value = -404.75
$ws.Range("N118").Value = -11186
$ws.Range("H121").Value = 1857.9354
$ws.Range("I121").Value = 368.33334
$ws.Range("J121").Value = 2017.5358
$ws.Range("K121").Value = 1105.00002
$ws.Range("L121").Value = 6052.607400000001
$ws.Range("M121").Value = 204.9999800000001
$ws.Range("N121").Value = -8672.607400000001
$ws.Range("H122").Value = 2588.8394
$ws.Range("J122").Value = 2930.3696
$ws.Range("L122").Value = 26373.3264
$ws.Range("N122").Value = -31273.3264
$ws.Range("H125").Value = 5138.3335
$ws.Range("I125").Value = 1950
$ws.Range("J125").Value = 8326.666999999999
$ws.Range("K125").Value = 5850
$ws.Range("L125").Value = 24980.001
$ws.Range("M125").Value = -930
$ws.Range("N125").Value = -34820.001
$ws.Range("H135").Value = 743402.1
$ws.Range("I135").Value = 660
$ws.Range("J135").Value = 1029072.1
$ws.Range("K135").Value = 5940
$ws.Range("L135").Value = 9261648.9
$ws.Range("M135").Value = -3405
$ws.Range("N135").Value = -9266718.9
$ws.Range("H141").Value = 8314.666999999999
$ws.Range("I141").Value = 8701.538
$ws.Range("K141").Value = 26104.614
$ws.Range("M141").Value = -20924.614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 5754.5
$ws.Range("J17").Value = 5754.5
$ws.Range("L17").Value = 5754.5
$ws.Range("N17").Value = -6090.5
$ws.Range("H132").Value = 3210.0417
$ws.Range("I132").Value = 1548.6
$ws.Range("J132").Value = 5979.1113
$ws.Range("K132").Value = 4645.799999999999
$ws.Range("L132").Value = 17937.3339
$ws.Range("M132").Value = -2115.799999999999
$ws.Range("N132").Value = -22997.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 30133
$ws.Range("J5").Value = 29999.5
$ws.Range("L5").Value = 29999.5
$ws.Range("N5").Value = -30225.5
$ws.Range("H40").Value = 6243.12
$ws.Range("I40").Value = 5321.727
$ws.Range("K40").Value = 5321.727
$ws.Range("M40").Value = -5185.727
$ws.Range("H46").Value = 1482.3077
$ws.Range("I46").Value = 842.6
$ws.Range("J46").Value = 2354.6365
$ws.Range("K46").Value = 842.6
$ws.Range("L46").Value = 2354.6365
$ws.Range("M46").Value = -654.6
$ws.Range("N46").Value = -2730.6365
$ws.Range("H115").Value = 38000
$ws.Range("J115").Value = 38000
$ws.Range("L115").Value = 38000
$ws.Range("N115").Value = -40350
$ws.Range("H122").Value = 4952.1055
$ws.Range("I122").Value = 4008.6924
$ws.Range("J122").Value = 6996.1665
$ws.Range("K122").Value = 12026.0772
$ws.Range("L122").Value = 20988.4995
$ws.Range("M122").Value = -9576.0772
$ws.Range("N122").Value = -25888.4995
$ws.Range("H132").Value = 5378.4443
$ws.Range("I132").Value = 4167.9165
$ws.Range("J132").Value = 7799.5
$ws.Range("K132").Value = 12503.7495
$ws.Range("L132").Value = 23398.5
$ws.Range("M132").Value = -9973.749500000002
$ws.Range("N132").Value = -28458.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H30").Value = 21562.25
$ws.Range("I30").Value = 16999
$ws.Range("J30").Value = 23083.334
$ws.Range("K30").Value = 16999
$ws.Range("L30").Value = 23083.334
$ws.Range("M30").Value = -16892
$ws.Range("N30").Value = -23297.334
$ws.Range("H41").Value = 19077.715
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 19077.715
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 19077.715
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -19857.715
$ws.Range("H132").Value = 20835976
$ws.Range("I132").Value = 1208.4286
$ws.Range("J132").Value = 37040790
$ws.Range("K132").Value = 3625.2858
$ws.Range("L132").Value = 111122370
$ws.Range("M132").Value = -1095.2858
$ws.Range("N132").Value = -111127430
$ws.Range("H138").Value = 42666.332
$ws.Range("J138").Value = 42666.332
$ws.Range("L138").Value = 42666.332
$ws.Range("N138").Value = -52946.332
